# Scheduled runner update: refresh market/profit figures across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 29915684
$ws.Range("I92").Value = 6173950.5
$ws.Range("K92").Value = 6173950.5
$ws.Range("M92").Value = -6172702.5
$ws.Range("H129").Value = 977.431
$ws.Range("I129").Value = 859.4
$ws.Range("J129").Value = 988.56604
$ws.Range("K129").Value = 2578.2
$ws.Range("L129").Value = 2965.69812
$ws.Range("M129").Value = 2421.8
$ws.Range("N129").Value = -12965.69812
$ws.Range("H132").Value = 1586.525
$ws.Range("I132").Value = 1354.1111
$ws.Range("J132").Value = 2069.2307
$ws.Range("K132").Value = 4062.3333
$ws.Range("L132").Value = 6207.6921
$ws.Range("M132").Value = -1532.3333
$ws.Range("N132").Value = -11267.6921
$ws.Range("H137").Value = 1240.1428
$ws.Range("I137").Value = 964.7805
$ws.Range("J137").Value = 1992.8
$ws.Range("K137").Value = 2894.3415
$ws.Range("L137").Value = 5978.4
$ws.Range("M137").Value = -344.3415
$ws.Range("N137").Value = -11078.4
$ws.Range("H138").Value = 2238.9011
$ws.Range("I138").Value = 956.2439
$ws.Range("J138").Value = 3290.68
$ws.Range("K138").Value = 2868.7317
$ws.Range("L138").Value = 9872.039999999999
$ws.Range("M138").Value = 2271.2683
$ws.Range("N138").Value = -20152.04
$ws.Range("H141").Value = 1564.6757
$ws.Range("I141").Value = 1070.1034
$ws.Range("J141").Value = 3357.5
$ws.Range("K141").Value = 3210.3102
$ws.Range("L141").Value = 10072.5
$ws.Range("M141").Value = 1969.6898
$ws.Range("N141").Value = -20432.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3758.63
$ws.Range("I32").Value = 2721.1584
$ws.Range("J32").Value = 8484.889
$ws.Range("K32").Value = 2721.1584
$ws.Range("L32").Value = 8484.889
$ws.Range("M32").Value = -2434.1584
$ws.Range("N32").Value = -9058.889
$ws.Range("H61").Value = 3069.173
$ws.Range("I61").Value = 3844.0278
$ws.Range("J61").Value = 1325.75
$ws.Range("K61").Value = 3844.0278
$ws.Range("L61").Value = 1325.75
$ws.Range("M61").Value = -3632.0278
$ws.Range("N61").Value = -1749.75
$ws.Range("H74").Value = 800.8182
$ws.Range("I74").Value = 660.3582
$ws.Range("J74").Value = 1248.9524
$ws.Range("K74").Value = 660.3582
$ws.Range("L74").Value = 1248.9524
$ws.Range("M74").Value = 213.6418
$ws.Range("N74").Value = -2996.9524
$ws.Range("H77").Value = 800.8182
$ws.Range("I77").Value = 660.3582
$ws.Range("J77").Value = 1248.9524
$ws.Range("K77").Value = 3301.791
$ws.Range("L77").Value = 6244.762
$ws.Range("M77").Value = 1066.209
$ws.Range("N77").Value = -14980.762
$ws.Range("H122").Value = 1832380.1
$ws.Range("I122").Value = 2331846.2
$ws.Range("J122").Value = 1004.6667
$ws.Range("K122").Value = 6995538.600000001
$ws.Range("L122").Value = 3014.0001
$ws.Range("M122").Value = -6993088.600000001
$ws.Range("N122").Value = -7914.0001
$ws.Range("H132").Value = 1889255.6
$ws.Range("I132").Value = 1687.9487
$ws.Range("J132").Value = 7147480
$ws.Range("K132").Value = 5063.8461
$ws.Range("L132").Value = 21442440
$ws.Range("M132").Value = -2533.8461
$ws.Range("N132").Value = -21447500
$ws.Range("H136").Value = 3069.173
$ws.Range("I136").Value = 3844.0278
$ws.Range("J136").Value = 1325.75
$ws.Range("K136").Value = 11532.0834
$ws.Range("L136").Value = 3977.25
$ws.Range("M136").Value = -8982.0834
$ws.Range("N136").Value = -9077.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 35275.5
$ws.Range("J126").Value = 35275.5
$ws.Range("L126").Value = 35275.5
$ws.Range("N126").Value = -45155.5
$ws.Range("H130").Value = 52779.668
$ws.Range("J130").Value = 52779.668
$ws.Range("L130").Value = 52779.668
$ws.Range("N130").Value = -62819.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 376.13333
$ws.Range("I22").Value = 304.63635
$ws.Range("J22").Value = 572.75
$ws.Range("K22").Value = 304.63635
$ws.Range("L22").Value = 572.75
$ws.Range("M22").Value = 45.36365000000001
$ws.Range("N22").Value = -1272.75
$ws.Range("H31").Value = 244500.64
$ws.Range("I31").Value = 1535.9348
$ws.Range("J31").Value = 865410.44
$ws.Range("K31").Value = 1535.9348
$ws.Range("L31").Value = 865410.44
$ws.Range("M31").Value = -1240.9348
$ws.Range("N31").Value = -866000.44
$ws.Range("H34").Value = 244500.64
$ws.Range("I34").Value = 1535.9348
$ws.Range("J34").Value = 865410.44
$ws.Range("K34").Value = 1535.9348
$ws.Range("L34").Value = 865410.44
$ws.Range("M34").Value = -1333.9348
$ws.Range("N34").Value = -865814.44
$ws.Range("H58").Value = 1612.0233
$ws.Range("I58").Value = 974.5417
$ws.Range("J58").Value = 2417.2632
$ws.Range("K58").Value = 974.5417
$ws.Range("L58").Value = 2417.2632
$ws.Range("M58").Value = -771.5417
$ws.Range("N58").Value = -2823.2632
$ws.Range("H132").Value = 1814.8209
$ws.Range("I132").Value = 1479.9375
$ws.Range("J132").Value = 2660.842
$ws.Range("K132").Value = 4439.8125
$ws.Range("L132").Value = 7982.526
$ws.Range("M132").Value = -1909.8125
$ws.Range("N132").Value = -13042.526
$ws.Range("H134").Value = 1755.9844
$ws.Range("I134").Value = 1979.8043
$ws.Range("J134").Value = 1184
$ws.Range("K134").Value = 5939.4129
$ws.Range("L134").Value = 3552
$ws.Range("M134").Value = -3404.4129
$ws.Range("N134").Value = -8622
$ws.Range("H136").Value = 1612.0233
$ws.Range("I136").Value = 974.5417
$ws.Range("J136").Value = 2417.2632
$ws.Range("K136").Value = 2923.6251
$ws.Range("L136").Value = 7251.7896
$ws.Range("M136").Value = -373.6251000000002
$ws.Range("N136").Value = -12351.7896
$ws.Range("H138").Value = 35846.668
$ws.Range("J138").Value = 35846.668
$ws.Range("L138").Value = 35846.668
$ws.Range("N138").Value = -46126.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 213316.92
$ws.Range("I113").Value = 563.7692
$ws.Range("J113").Value = 476725.56
$ws.Range("K113").Value = 1691.3076
$ws.Range("L113").Value = 1430176.68
$ws.Range("M113").Value = 478.6924000000001
$ws.Range("N113").Value = -1434516.68
$ws.Range("H131").Value = 3125907.5
$ws.Range("I131").Value = 14286067
$ws.Range("J131").Value = 1062.84
$ws.Range("K131").Value = 42858201
$ws.Range("L131").Value = 3188.52
$ws.Range("M131").Value = -42853161
$ws.Range("N131").Value = -13268.52

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4223.143
$ws.Range("I126").Value = 6439.2383
$ws.Range("J126").Value = 2561.0715
$ws.Range("K126").Value = 19317.7149
$ws.Range("L126").Value = 7683.2145
$ws.Range("M126").Value = -16847.7149
$ws.Range("N126").Value = -12623.2145
$ws.Range("H132").Value = 1989.4117
$ws.Range("I132").Value = 1612.0741
$ws.Range("J132").Value = 2413.9167
$ws.Range("K132").Value = 4836.2223
$ws.Range("L132").Value = 7241.750100000001
$ws.Range("M132").Value = -2306.2223
$ws.Range("N132").Value = -12301.7501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9377974
$ws.Range("I132").Value = 12147770
$ws.Range("J132").Value = 3277.1538
$ws.Range("K132").Value = 36443310
$ws.Range("L132").Value = 9831.4614
$ws.Range("M132").Value = -36440780
$ws.Range("N132").Value = -14891.4614
$ws.Range("H136").Value = 8301.535
$ws.Range("I136").Value = 5756.355
$ws.Range("J136").Value = 14876.583
$ws.Range("K136").Value = 17269.065
$ws.Range("L136").Value = 44629.749
$ws.Range("M136").Value = -14719.065
$ws.Range("N136").Value = -49729.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 38666.668
$ws.Range("J63").Value = 38666.668
$ws.Range("L63").Value = 38666.668
$ws.Range("N63").Value = -39914.668
$ws.Range("H66").Value = 38666.668
$ws.Range("J66").Value = 38666.668
$ws.Range("L66").Value = 116000.004
$ws.Range("N66").Value = -122240.004
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 1812.7632
$ws.Range("I122").Value = 1765
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 5295
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -2845
$ws.Range("N122").Value = -10799.9998
$ws.Range("H132").Value = 13233.35
$ws.Range("I132").Value = 16425.016
$ws.Range("J132").Value = 1405.4117
$ws.Range("K132").Value = 49275.048
$ws.Range("L132").Value = 4216.2351
$ws.Range("M132").Value = -46745.048
$ws.Range("N132").Value = -9276.2351
$ws.Range("H136").Value = 6495609
$ws.Range("I136").Value = 2052.3774
$ws.Range("K136").Value = 6157.1322
$ws.Range("M136").Value = -3607.1322
